$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a run of year labels ("2016".."2050") as TEXT (not numbers)
# into a 1-row helper range, then copy/paste-values it into the real
# destination range so the destination cells end up as plain shared-string
# text cells (matching how the original workbook stores "1500".."2015" etc.)
# instead of picking up an explicit NumberFormat style.
# ---------------------------------------------------------------------------
function Write-YearHeaders {
    param($ws, $helperRow, $startCol, $destRow, $destColStart, $firstYear, $lastYear)

    $n = $lastYear - $firstYear + 1
    for ($i = 0; $i -lt $n; $i++) {
        $h = $ws.Cells.Item($helperRow, $startCol + $i)
        $h.NumberFormat = "@"
        $h.Value = [string]($firstYear + $i)
        $h.ClearFormats()
    }

    $helperRange = $ws.Range($ws.Cells.Item($helperRow, $startCol), $ws.Cells.Item($helperRow, $startCol + $n - 1))
    $helperRange.Copy()
    $destRange = $ws.Range($ws.Cells.Item($destRow, $destColStart), $ws.Cells.Item($destRow, $destColStart + $n - 1))
    $destRange.PasteSpecial(-4163)   # xlPasteValues
    $helperRange.EntireRow.Delete()
}

# ---------------------------------------------------------------------------
# 1) "Data Clio Infra Format" sheet: append year columns 2016-2050 (wide
#    format). They go right after the existing last year column (2015 =
#    column TD), so 2016 starts at column TE (index 525). Row 2 (the Ceylon
#    data row) has no data for these new years, so it stays blank.
# ---------------------------------------------------------------------------
$wsWide = $wb.Worksheets.Item("Data Clio Infra Format")
Write-YearHeaders $wsWide 100 525 1 525 2016 2050

# ---------------------------------------------------------------------------
# 2) "Data Long Format" sheet: insert 35 new columns before the "year"
#    column (E), pushing "year"/"value" from E/F to AN/AO, and stamp the new
#    columns' header row with the same 2016-2050 labels. The inserted data
#    cells (rows 2-14) stay blank - there's no long-format data yet for
#    those years.
# ---------------------------------------------------------------------------
$wsLong = $wb.Worksheets.Item("Data Long Format")
$wsLong.Columns("E:AM").Insert()
Write-YearHeaders $wsLong 100 5 1 5 2016 2050

# ---------------------------------------------------------------------------
# 3) "Metadata" sheet: citation text gained an Oxford comma.
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("C3").Value = "Zwart, Pim de, Bas van Leeuwen, and Jieli van Leeuwen-Li (2015). Labourers Real Wage. http://hdl.handle.net/10622/QK8VRF, accessed via the Clio Infra website."
